$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2052238805970149
$ws.Range("C2").Value = 0.5447761194029851
$ws.Range("J2").Value = 0.01492537313432836
$ws.Range("P2").Value = 0.1380597014925373
$ws.Range("S2").Value = 0.09701492537313433
$ws.Range("C3").Value = 0.006666666666666667
$ws.Range("J3").Value = 0.02
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.1733333333333333
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6470588235294118
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.04484304932735426
$ws.Range("F6").Value = 0.1031390134529148
$ws.Range("J6").Value = 0.2556053811659193
$ws.Range("O6").Value = 0.02242152466367713
$ws.Range("Q6").Value = 0.1524663677130045
$ws.Range("R6").Value = 0.07174887892376682
$ws.Range("S6").Value = 0.3497757847533632
$ws.Range("B7").Value = 0.1325301204819277
$ws.Range("D7").Value = 0.02409638554216868
$ws.Range("F7").Value = 0.0783132530120482
$ws.Range("J7").Value = 0.1024096385542169
$ws.Range("O7").Value = 0.006024096385542169
$ws.Range("Q7").Value = 0.1807228915662651
$ws.Range("R7").Value = 0.0783132530120482
$ws.Range("S7").Value = 0.3975903614457831
$ws.Range("B8").Value = 0.09510086455331412
$ws.Range("D8").Value = 0.01440922190201729
$ws.Range("F8").Value = 0.06628242074927954
$ws.Range("J8").Value = 0.09221902017291066
$ws.Range("O8").Value = 0.01729106628242075
$ws.Range("Q8").Value = 0.1152737752161383
$ws.Range("R8").Value = 0.1671469740634006
$ws.Range("S8").Value = 0.4322766570605187
$ws.Range("B9").Value = 0.1012658227848101
$ws.Range("D9").Value = 0.01265822784810127
$ws.Range("F9").Value = 0.109704641350211
$ws.Range("J9").Value = 0.09282700421940929
$ws.Range("O9").Value = 0.008438818565400843
$ws.Range("Q9").Value = 0.1561181434599156
$ws.Range("R9").Value = 0.1350210970464135
$ws.Range("S9").Value = 0.3839662447257384
$ws.Range("B10").Value = 0.09975669099756691
$ws.Range("D10").Value = 0.0178426601784266
$ws.Range("E10").Value = 0.0008110300081103001
$ws.Range("F10").Value = 0.0575831305758313
$ws.Range("J10").Value = 0.1240875912408759
$ws.Range("O10").Value = 0.0202757502027575
$ws.Range("Q10").Value = 0.1987023519870235
$ws.Range("R10").Value = 0.1200324412003244
$ws.Range("S10").Value = 0.3609083536090835
$ws.Range("G11").Value = 0.138996138996139
$ws.Range("J11").Value = 0.09266409266409266
$ws.Range("K11").Value = 0.1969111969111969
$ws.Range("L11").Value = 0.5521235521235521
$ws.Range("S11").Value = 0.0193050193050193
$ws.Range("G12").Value = 0.7569444444444444
$ws.Range("J12").Value = 0.1805555555555556
$ws.Range("K12").Value = 0.01388888888888889
$ws.Range("L12").Value = 0.02777777777777778
$ws.Range("S12").Value = 0.02083333333333333
$ws.Range("G13").Value = 0.6444444444444445
$ws.Range("J13").Value = 0.2888888888888889
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.005154639175257732
$ws.Range("H15").Value = 0.1391752577319588
$ws.Range("I15").Value = 0.09793814432989691
$ws.Range("J15").Value = 0.3762886597938144
$ws.Range("K15").Value = 0.07216494845360824
$ws.Range("M15").Value = 0.0154639175257732
$ws.Range("O15").Value = 0.09278350515463918
$ws.Range("S15").Value = 0.2010309278350516
$ws.Range("F16").Value = 0.02325581395348837
$ws.Range("H16").Value = 0.1511627906976744
$ws.Range("I16").Value = 0.0872093023255814
$ws.Range("J16").Value = 0.4127906976744186
$ws.Range("K16").Value = 0.1220930232558139
$ws.Range("M16").Value = 0.02325581395348837
$ws.Range("N16").Value = 0.005813953488372093
$ws.Range("O16").Value = 0.05813953488372093
$ws.Range("S16").Value = 0.1162790697674419
$ws.Range("F17").Value = 0.01044386422976501
$ws.Range("H17").Value = 0.1331592689295039
$ws.Range("I17").Value = 0.1409921671018277
$ws.Range("J17").Value = 0.4308093994778068
$ws.Range("K17").Value = 0.08616187989556136
$ws.Range("M17").Value = 0.01827676240208877
$ws.Range("O17").Value = 0.04960835509138381
$ws.Range("S17").Value = 0.1305483028720627
$ws.Range("F18").Value = 0.01123595505617977
$ws.Range("H18").Value = 0.1460674157303371
$ws.Range("I18").Value = 0.0898876404494382
$ws.Range("J18").Value = 0.4419475655430712
$ws.Range("K18").Value = 0.1048689138576779
$ws.Range("M18").Value = 0.02247191011235955
$ws.Range("O18").Value = 0.06367041198501873
$ws.Range("S18").Value = 0.1198501872659176
$ws.Range("F19").Value = 0.01862828111769687
$ws.Range("H19").Value = 0.1727349703640982
$ws.Range("I19").Value = 0.104995766299746
$ws.Range("J19").Value = 0.3979678238780694
$ws.Range("K19").Value = 0.08975444538526672
$ws.Range("M19").Value = 0.02116850127011008
$ws.Range("N19").Value = 0.000846740050804403
$ws.Range("O19").Value = 0.06350550381033022
$ws.Range("S19").Value = 0.1303979678238781
